$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 232, pushing the existing rows 232-272 down to 233-273.
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row with the new Jengibre price-record data
# (same template as the surrounding rows, new date/volume/price figures).
$ws.Cells.Item(232, 1).Value = 10
$ws.Cells.Item(232, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(232, 3).Value = "La Araucanía"
$ws.Cells.Item(232, 4).Value = 45015
$ws.Cells.Item(232, 5).Value = 9
$ws.Cells.Item(232, 6).Value = 100114007
$ws.Cells.Item(232, 7).Value = "Jengibre"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 80
$ws.Cells.Item(232, 11).Value = 25000
$ws.Cells.Item(232, 12).Value = 25000
$ws.Cells.Item(232, 13).Value = 25000
$ws.Cells.Item(232, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(232, 15).Value = "Perú"
$ws.Cells.Item(232, 16).Value = 1923
$ws.Cells.Item(232, 17).Value = 13
$ws.Cells.Item(232, 18).Value = "Hortaliza"

Write-Output "inserted row 232; sheet now $($ws.Range("A1").CurrentRegion.Rows.Count) rows used"
